$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 470.01666  # H17: 503.7414 -> 470.01666
$ws.Cells.Item(17, 10).Value = 480.3684  # J17: 516.3090999999999 -> 480.3684
$ws.Cells.Item(17, 12).Value = 1441.1052  # L17: 1548.9273 -> 1441.1052
$ws.Cells.Item(17, 14).Value = -1777.1052  # N17: -1884.9273 -> -1777.1052
$ws.Cells.Item(106, 8).Value = 3828.2222  # H106: 4187.375 -> 3828.2222
$ws.Cells.Item(106, 9).Value = 1151  # I106: 1249 -> 1151
$ws.Cells.Item(106, 11).Value = 1151  # K106: 1249 -> 1151
$ws.Cells.Item(106, 13).Value = -520  # M106: -618 -> -520
$ws.Cells.Item(111, 8).Value = 3003.75  # H111: 3022.0833 -> 3003.75
$ws.Cells.Item(111, 9).Value = 2913.2727  # I111: 3126.6 -> 2913.2727
$ws.Cells.Item(111, 10).Value = 3999  # J111: 2499.5 -> 3999
$ws.Cells.Item(111, 11).Value = 8739.8181  # K111: 9379.799999999999 -> 8739.8181
$ws.Cells.Item(111, 12).Value = 11997  # L111: 7498.5 -> 11997
$ws.Cells.Item(111, 13).Value = -5672.8181  # M111: -6312.799999999999 -> -5672.8181
$ws.Cells.Item(111, 14).Value = -18131  # N111: -13632.5 -> -18131
$ws.Cells.Item(113, 8).Value = 3874.75  # H113: 4250 -> 3874.75
$ws.Cells.Item(113, 9).Value = 2833  # I113: 2875 -> 2833
$ws.Cells.Item(113, 11).Value = 2833  # K113: 2875 -> 2833
$ws.Cells.Item(113, 13).Value = 421  # M113: 379 -> 421
$ws.Cells.Item(116, 8).Value = 5948.533  # H116: 6352.143 -> 5948.533
$ws.Cells.Item(116, 9).Value = 5699.875  # I116: 6100 -> 5699.875
$ws.Cells.Item(116, 10).Value = 6232.7144  # J116: 6604.2856 -> 6232.7144
$ws.Cells.Item(116, 11).Value = 5699.875  # K116: 6100 -> 5699.875
$ws.Cells.Item(116, 12).Value = 6232.7144  # L116: 6604.2856 -> 6232.7144
$ws.Cells.Item(116, 13).Value = -2257.875  # M116: -2658 -> -2257.875
$ws.Cells.Item(116, 14).Value = -13116.7144  # N116: -13488.2856 -> -13116.7144
$ws.Cells.Item(137, 8).Value = 3453.6758  # H137: 3259.1462 -> 3453.6758
$ws.Cells.Item(137, 9).Value = 1724.7916  # I137: 1660.5186 -> 1724.7916
$ws.Cells.Item(137, 10).Value = 6645.4614  # J137: 6342.2144 -> 6645.4614
$ws.Cells.Item(137, 11).Value = 5174.3748  # K137: 4981.5558 -> 5174.3748
$ws.Cells.Item(137, 12).Value = 19936.3842  # L137: 19026.6432 -> 19936.3842
$ws.Cells.Item(137, 13).Value = -2624.3748  # M137: -2431.5558 -> -2624.3748
$ws.Cells.Item(137, 14).Value = -25036.3842  # N137: -24126.6432 -> -25036.3842
$ws.Cells.Item(138, 8).Value = 2017.1719  # H138: 2038.9048 -> 2017.1719
$ws.Cells.Item(138, 9).Value = 868.2222  # I138: 881.05884 -> 868.2222
$ws.Cells.Item(138, 10).Value = 2466.761  # J138: 2466.8044 -> 2466.761
$ws.Cells.Item(138, 11).Value = 2604.6666  # K138: 2643.17652 -> 2604.6666
$ws.Cells.Item(138, 12).Value = 7400.282999999999  # L138: 7400.4132 -> 7400.282999999999
$ws.Cells.Item(138, 13).Value = 2535.3334  # M138: 2496.82348 -> 2535.3334
$ws.Cells.Item(138, 14).Value = -17680.283  # N138: -17680.4132 -> -17680.283

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13160405  # H32: 13891541 -> 13160405
$ws.Cells.Item(32, 9).Value = 14707414  # I32: 15626627 -> 14707414
$ws.Cells.Item(32, 10).Value = 10824.25  # J32: 10849.5 -> 10824.25
$ws.Cells.Item(32, 11).Value = 14707414  # K32: 15626627 -> 14707414
$ws.Cells.Item(32, 12).Value = 10824.25  # L32: 10849.5 -> 10824.25
$ws.Cells.Item(32, 13).Value = -14707127  # M32: -15626340 -> -14707127
$ws.Cells.Item(32, 14).Value = -11398.25  # N32: -11423.5 -> -11398.25
$ws.Cells.Item(61, 8).Value = 13547208  # H61: 14321305 -> 13547208
$ws.Cells.Item(61, 9).Value = 16671667  # I61: 17862464 -> 16671667
$ws.Cells.Item(61, 11).Value = 16671667  # K61: 17862464 -> 16671667
$ws.Cells.Item(61, 13).Value = -16671455  # M61: -17862252 -> -16671455
$ws.Cells.Item(122, 8).Value = 2593.2666  # H122: 2883.25 -> 2593.2666
$ws.Cells.Item(122, 9).Value = 1825  # I122: 3000 -> 1825
$ws.Cells.Item(122, 11).Value = 5475  # K122: 9000 -> 5475
$ws.Cells.Item(122, 13).Value = -3025  # M122: -6550 -> -3025
$ws.Cells.Item(132, 8).Value = 11217  # H132: 11792.3125 -> 11217
$ws.Cells.Item(132, 9).Value = 5334.3335  # I132: 5749.625 -> 5334.3335
$ws.Cells.Item(132, 11).Value = 16003.0005  # K132: 17248.875 -> 16003.0005
$ws.Cells.Item(132, 13).Value = -13473.0005  # M132: -14718.875 -> -13473.0005
$ws.Cells.Item(136, 8).Value = 13547208  # H136: 14321305 -> 13547208
$ws.Cells.Item(136, 9).Value = 16671667  # I136: 17862464 -> 16671667
$ws.Cells.Item(136, 11).Value = 50015001  # K136: 53587392 -> 50015001
$ws.Cells.Item(136, 13).Value = -50012451  # M136: -53584842 -> -50012451

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(69, 8).Value = 60000  # H69: 0 -> 60000
$ws.Cells.Item(69, 10).Value = 60000  # J69: 0 -> 60000
$ws.Cells.Item(69, 12).Value = 60000  # L69: 0 -> 60000
$ws.Cells.Item(69, 14).Value = -61622  # N69: None -> -61622
$ws.Cells.Item(72, 8).Value = 60000  # H72: 0 -> 60000
$ws.Cells.Item(72, 10).Value = 60000  # J72: 0 -> 60000
$ws.Cells.Item(72, 12).Value = 180000  # L72: 0 -> 180000
$ws.Cells.Item(72, 14).Value = -188112  # N72: None -> -188112
$ws.Cells.Item(99, 8).Value = 2138.8572  # H99: 2159.6 -> 2138.8572
$ws.Cells.Item(99, 9).Value = 1678.6666  # I99: 1738 -> 1678.6666
$ws.Cells.Item(99, 11).Value = 1678.6666  # K99: 1738 -> 1678.6666
$ws.Cells.Item(99, 13).Value = -180.6666  # M99: -240 -> -180.6666
$ws.Cells.Item(103, 8).Value = 30330.6  # H103: 14249.5 -> 30330.6
$ws.Cells.Item(103, 10).Value = 30330.6  # J103: 14249.5 -> 30330.6
$ws.Cells.Item(103, 12).Value = 30330.6  # L103: 14249.5 -> 30330.6
$ws.Cells.Item(103, 14).Value = -32674.6  # N103: -16593.5 -> -32674.6
$ws.Cells.Item(106, 8).Value = 63316.332  # H106: 63296.332 -> 63316.332
$ws.Cells.Item(106, 10).Value = 63316.332  # J106: 63296.332 -> 63316.332
$ws.Cells.Item(106, 12).Value = 63316.332  # L106: 63296.332 -> 63316.332
$ws.Cells.Item(106, 14).Value = -65840.33199999999  # N106: -65820.33199999999 -> -65840.33199999999
$ws.Cells.Item(109, 8).Value = 73684  # H109: 79831.5 -> 73684
$ws.Cells.Item(109, 10).Value = 73684  # J109: 79831.5 -> 73684
$ws.Cells.Item(109, 12).Value = 73684  # L109: 79831.5 -> 73684
$ws.Cells.Item(109, 14).Value = -76458  # N109: -82605.5 -> -76458
$ws.Cells.Item(134, 8).Value = 181741.17  # H134: 154942.28 -> 181741.17
$ws.Cells.Item(134, 9).Value = 10209.5  # I134: 4856 -> 10209.5
$ws.Cells.Item(134, 11).Value = 30628.5  # K134: 14568 -> 30628.5
$ws.Cells.Item(134, 13).Value = -28093.5  # M134: -12033 -> -28093.5
$ws.Cells.Item(140, 8).Value = 0  # H140: 100000 -> 0
$ws.Cells.Item(140, 10).Value = 0  # J140: 100000 -> 0
$ws.Cells.Item(140, 12).Value = 0  # L140: 100000 -> 0
$ws.Cells.Item(140, 14).ClearContents()  # N140: -110360 -> (blank)
$ws.Cells.Item(141, 8).Value = 0  # H141: 75000 -> 0
$ws.Cells.Item(141, 10).Value = 0  # J141: 75000 -> 0
$ws.Cells.Item(141, 12).Value = 0  # L141: 75000 -> 0
$ws.Cells.Item(141, 14).ClearContents()  # N141: -85360 -> (blank)

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(44, 8).Value = 64  # H44: 0 -> 64
$ws.Cells.Item(44, 9).Value = 64  # I44: 0 -> 64
$ws.Cells.Item(44, 11).Value = 64  # K44: 0 -> 64
$ws.Cells.Item(44, 13).Value = 378  # M44: None -> 378
$ws.Cells.Item(86, 8).Value = 2466.6667  # H86: 2500 -> 2466.6667
$ws.Cells.Item(86, 10).Value = 2950  # J86: 3000 -> 2950
$ws.Cells.Item(86, 12).Value = 2950  # L86: 3000 -> 2950
$ws.Cells.Item(86, 14).Value = -5196  # N86: -5246 -> -5196
$ws.Cells.Item(89, 8).Value = 2466.6667  # H89: 2500 -> 2466.6667
$ws.Cells.Item(89, 10).Value = 2950  # J89: 3000 -> 2950
$ws.Cells.Item(89, 12).Value = 14750  # L89: 15000 -> 14750
$ws.Cells.Item(89, 14).Value = -25982  # N89: -26232 -> -25982
$ws.Cells.Item(99, 8).Value = 2535.2856  # H99: 2753.5386 -> 2535.2856
$ws.Cells.Item(99, 9).Value = 2246  # I99: 2531.8572 -> 2246
$ws.Cells.Item(99, 10).Value = 2824.5715  # J99: 3012.1667 -> 2824.5715
$ws.Cells.Item(99, 11).Value = 2246  # K99: 2531.8572 -> 2246
$ws.Cells.Item(99, 12).Value = 2824.5715  # L99: 3012.1667 -> 2824.5715
$ws.Cells.Item(99, 13).Value = -748  # M99: -1033.8572 -> -748
$ws.Cells.Item(99, 14).Value = -5820.5715  # N99: -6008.1667 -> -5820.5715
$ws.Cells.Item(126, 8).Value = 2535.2856  # H126: 2753.5386 -> 2535.2856
$ws.Cells.Item(126, 9).Value = 2246  # I126: 2531.8572 -> 2246
$ws.Cells.Item(126, 10).Value = 2824.5715  # J126: 3012.1667 -> 2824.5715
$ws.Cells.Item(126, 11).Value = 6738  # K126: 7595.571599999999 -> 6738
$ws.Cells.Item(126, 12).Value = 8473.7145  # L126: 9036.500100000001 -> 8473.7145
$ws.Cells.Item(126, 13).Value = -4268  # M126: -5125.571599999999 -> -4268
$ws.Cells.Item(126, 14).Value = -13413.7145  # N126: -13976.5001 -> -13413.7145
$ws.Cells.Item(132, 8).Value = 2269.9  # H132: 2411 -> 2269.9
$ws.Cells.Item(132, 9).Value = 2269.9  # I132: 2411 -> 2269.9
$ws.Cells.Item(132, 11).Value = 6809.700000000001  # K132: 7233 -> 6809.700000000001
$ws.Cells.Item(132, 13).Value = -4279.700000000001  # M132: -4703 -> -4279.700000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 2169.5  # H122: 2065.1428 -> 2169.5
$ws.Cells.Item(122, 10).Value = 2443.7144  # J122: 2367.6667 -> 2443.7144
$ws.Cells.Item(122, 12).Value = 21993.4296  # L122: 21309.0003 -> 21993.4296
$ws.Cells.Item(122, 14).Value = -26893.4296  # N122: -26209.0003 -> -26893.4296
$ws.Cells.Item(125, 8).Value = 13433.286  # H125: 13147.429 -> 13433.286
$ws.Cells.Item(125, 9).Value = 0  # I125: 1999 -> 0
$ws.Cells.Item(125, 10).Value = 13433.286  # J125: 15005.5 -> 13433.286
$ws.Cells.Item(125, 11).Value = 0  # K125: 5997 -> 0
$ws.Cells.Item(125, 12).Value = 40299.858  # L125: 45016.5 -> 40299.858
$ws.Cells.Item(125, 13).ClearContents()  # M125: -1077 -> (blank)
$ws.Cells.Item(125, 14).Value = -50139.858  # N125: -54856.5 -> -50139.858
$ws.Cells.Item(139, 8).Value = 2675  # H139: 2631.5789 -> 2675
$ws.Cells.Item(139, 9).Value = 3000  # I139: 2500 -> 3000
$ws.Cells.Item(139, 11).Value = 9000  # K139: 7500 -> 9000
$ws.Cells.Item(139, 13).Value = -3860  # M139: -2360 -> -3860

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 3952.682  # H70: 4375.9414 -> 3952.682
$ws.Cells.Item(70, 9).Value = 4031.6924  # I70: 4525.4443 -> 4031.6924
$ws.Cells.Item(70, 10).Value = 3838.5557  # J70: 4207.75 -> 3838.5557
$ws.Cells.Item(70, 11).Value = 4031.6924  # K70: 4525.4443 -> 4031.6924
$ws.Cells.Item(70, 12).Value = 3838.5557  # L70: 4207.75 -> 3838.5557
$ws.Cells.Item(70, 13).Value = -3761.6924  # M70: -4255.4443 -> -3761.6924
$ws.Cells.Item(70, 14).Value = -4378.5557  # N70: -4747.75 -> -4378.5557
$ws.Cells.Item(73, 8).Value = 3952.682  # H73: 4375.9414 -> 3952.682
$ws.Cells.Item(73, 9).Value = 4031.6924  # I73: 4525.4443 -> 4031.6924
$ws.Cells.Item(73, 10).Value = 3838.5557  # J73: 4207.75 -> 3838.5557
$ws.Cells.Item(73, 11).Value = 4031.6924  # K73: 4525.4443 -> 4031.6924
$ws.Cells.Item(73, 12).Value = 3838.5557  # L73: 4207.75 -> 3838.5557
$ws.Cells.Item(73, 13).Value = -3095.6924  # M73: -3589.4443 -> -3095.6924
$ws.Cells.Item(73, 14).Value = -5710.5557  # N73: -6079.75 -> -5710.5557
$ws.Cells.Item(102, 8).Value = 5584.727  # H102: 5322.091 -> 5584.727
$ws.Cells.Item(102, 9).Value = 3554.125  # I102: 3282.5557 -> 3554.125
$ws.Cells.Item(102, 10).Value = 10999.667  # J102: 14500 -> 10999.667
$ws.Cells.Item(102, 11).Value = 3554.125  # K102: 3282.5557 -> 3554.125
$ws.Cells.Item(102, 12).Value = 10999.667  # L102: 14500 -> 10999.667
$ws.Cells.Item(102, 13).Value = -1932.125  # M102: -1660.5557 -> -1932.125
$ws.Cells.Item(102, 14).Value = -14243.667  # N102: -17744 -> -14243.667
$ws.Cells.Item(122, 8).Value = 2133.8  # H122: 1956.3334 -> 2133.8
$ws.Cells.Item(122, 9).Value = 2167.5  # I122: 1935 -> 2167.5
$ws.Cells.Item(122, 11).Value = 6502.5  # K122: 5805 -> 6502.5
$ws.Cells.Item(122, 13).Value = -4052.5  # M122: -3355 -> -4052.5
$ws.Cells.Item(132, 8).Value = 52647308  # H132: 58840530 -> 52647308
$ws.Cells.Item(132, 9).Value = 66672590  # I132: 76929150 -> 66672590
$ws.Cells.Item(132, 11).Value = 200017770  # K132: 230787450 -> 200017770
$ws.Cells.Item(132, 13).Value = -200015240  # M132: -230784920 -> -200015240
$ws.Cells.Item(136, 8).Value = 39900  # H136: 40217.332 -> 39900
$ws.Cells.Item(136, 10).Value = 39900  # J136: 40217.332 -> 39900
$ws.Cells.Item(136, 12).Value = 119700  # L136: 120651.996 -> 119700
$ws.Cells.Item(136, 14).Value = -124800  # N136: -125751.996 -> -124800

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 18276576  # H7: 18276602 -> 18276576
$ws.Cells.Item(7, 9).Value = 40001700  # I7: 50001876 -> 40001700
$ws.Cells.Item(7, 10).Value = 172304.17  # J7: 147875 -> 172304.17
$ws.Cells.Item(7, 11).Value = 40001700  # K7: 50001876 -> 40001700
$ws.Cells.Item(7, 12).Value = 172304.17  # L7: 147875 -> 172304.17
$ws.Cells.Item(7, 13).Value = -40001588  # M7: -50001764 -> -40001588
$ws.Cells.Item(7, 14).Value = -172528.17  # N7: -148099 -> -172528.17
$ws.Cells.Item(22, 8).Value = 2378.4  # H22: 2950.75 -> 2378.4
$ws.Cells.Item(22, 9).Value = 2899.25  # I22: 2950.75 -> 2899.25
$ws.Cells.Item(22, 10).Value = 295  # J22: 0 -> 295
$ws.Cells.Item(22, 11).Value = 2899.25  # K22: 2950.75 -> 2899.25
$ws.Cells.Item(22, 12).Value = 295  # L22: 0 -> 295
$ws.Cells.Item(22, 13).Value = -2604.25  # M22: -2655.75 -> -2604.25
$ws.Cells.Item(22, 14).Value = -885  # N22: None -> -885
$ws.Cells.Item(27, 8).Value = 2378.4  # H27: 2950.75 -> 2378.4
$ws.Cells.Item(27, 9).Value = 2899.25  # I27: 2950.75 -> 2899.25
$ws.Cells.Item(27, 10).Value = 295  # J27: 0 -> 295
$ws.Cells.Item(27, 11).Value = 2899.25  # K27: 2950.75 -> 2899.25
$ws.Cells.Item(27, 12).Value = 295  # L27: 0 -> 295
$ws.Cells.Item(27, 13).Value = -2792.25  # M27: -2843.75 -> -2792.25
$ws.Cells.Item(27, 14).Value = -509  # N27: None -> -509
$ws.Cells.Item(40, 8).Value = 3472.2666  # H40: 3443.7932 -> 3472.2666
$ws.Cells.Item(40, 10).Value = 5362.125  # J40: 5514.143 -> 5362.125
$ws.Cells.Item(40, 12).Value = 5362.125  # L40: 5514.143 -> 5362.125
$ws.Cells.Item(40, 14).Value = -5634.125  # N40: -5786.143 -> -5634.125
$ws.Cells.Item(61, 8).Value = 1045.7778  # H61: 797.625 -> 1045.7778
$ws.Cells.Item(61, 9).Value = 1068.3334  # I61: 750.7692 -> 1068.3334
$ws.Cells.Item(61, 11).Value = 1068.3334  # K61: 750.7692 -> 1068.3334
$ws.Cells.Item(61, 13).Value = -866.3334  # M61: -548.7692 -> -866.3334
$ws.Cells.Item(113, 8).Value = 1045.7778  # H113: 797.625 -> 1045.7778
$ws.Cells.Item(113, 9).Value = 1068.3334  # I113: 750.7692 -> 1068.3334
$ws.Cells.Item(113, 11).Value = 1068.3334  # K113: 750.7692 -> 1068.3334
$ws.Cells.Item(113, 13).Value = 1101.6666  # M113: 1419.2308 -> 1101.6666
$ws.Cells.Item(126, 8).Value = 18276576  # H126: 18276602 -> 18276576
$ws.Cells.Item(126, 9).Value = 40001700  # I126: 50001876 -> 40001700
$ws.Cells.Item(126, 10).Value = 172304.17  # J126: 147875 -> 172304.17
$ws.Cells.Item(126, 11).Value = 120005100  # K126: 150005628 -> 120005100
$ws.Cells.Item(126, 12).Value = 516912.51  # L126: 443625 -> 516912.51
$ws.Cells.Item(126, 13).Value = -120002630  # M126: -150003158 -> -120002630
$ws.Cells.Item(126, 14).Value = -521852.51  # N126: -448565 -> -521852.51
$ws.Cells.Item(136, 8).Value = 33844.08  # H136: 34859.656 -> 33844.08
$ws.Cells.Item(136, 9).Value = 4551.48  # I136: 4699.4585 -> 4551.48
$ws.Cells.Item(136, 10).Value = 86152.28999999999  # J136: 86562.86 -> 86152.28999999999
$ws.Cells.Item(136, 11).Value = 13654.44  # K136: 14098.3755 -> 13654.44
$ws.Cells.Item(136, 12).Value = 258456.87  # L136: 259688.58 -> 258456.87
$ws.Cells.Item(136, 13).Value = -11104.44  # M136: -11548.3755 -> -11104.44
$ws.Cells.Item(136, 14).Value = -263556.87  # N136: -264788.58 -> -263556.87

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 0  # H32: 17000 -> 0
$ws.Cells.Item(32, 9).Value = 0  # I32: 17000 -> 0
$ws.Cells.Item(32, 11).Value = 0  # K32: 17000 -> 0
$ws.Cells.Item(32, 13).ClearContents()  # M32: -16683 -> (blank)
$ws.Cells.Item(45, 8).Value = 8563  # H45: 8525.333000000001 -> 8563
$ws.Cells.Item(45, 10).Value = 8563  # J45: 8525.333000000001 -> 8563
$ws.Cells.Item(45, 12).Value = 8563  # L45: 8525.333000000001 -> 8563
$ws.Cells.Item(45, 14).Value = -9545  # N45: -9507.333000000001 -> -9545
$ws.Cells.Item(46, 8).Value = 49996.5  # H46: 53331 -> 49996.5
$ws.Cells.Item(46, 10).Value = 49996.5  # J46: 53331 -> 49996.5
$ws.Cells.Item(46, 12).Value = 49996.5  # L46: 53331 -> 49996.5
$ws.Cells.Item(46, 14).Value = -50458.5  # N46: -53793 -> -50458.5
$ws.Cells.Item(97, 8).Value = 112000  # H97: 19786 -> 112000
$ws.Cells.Item(97, 10).Value = 112000  # J97: 19786 -> 112000
$ws.Cells.Item(97, 12).Value = 112000  # L97: 19786 -> 112000
$ws.Cells.Item(97, 14).Value = -113982  # N97: -21768 -> -113982
$ws.Cells.Item(122, 8).Value = 8749.294  # H122: 9762.6 -> 8749.294
$ws.Cells.Item(122, 9).Value = 3212.25  # I122: 3899.8333 -> 3212.25
$ws.Cells.Item(122, 11).Value = 9636.75  # K122: 11699.4999 -> 9636.75
$ws.Cells.Item(122, 13).Value = -7186.75  # M122: -9249.499899999999 -> -7186.75
$ws.Cells.Item(131, 8).Value = 78904  # H131: 0 -> 78904
$ws.Cells.Item(131, 10).Value = 78904  # J131: 0 -> 78904
$ws.Cells.Item(131, 12).Value = 78904  # L131: 0 -> 78904
$ws.Cells.Item(131, 14).Value = -88984  # N131: None -> -88984
$ws.Cells.Item(132, 8).Value = 18049.334  # H132: 14702.733 -> 18049.334
$ws.Cells.Item(132, 9).Value = 1765.3334  # I132: 1615.6666 -> 1765.3334
$ws.Cells.Item(132, 11).Value = 5296.0002  # K132: 4846.9998 -> 5296.0002
$ws.Cells.Item(132, 13).Value = -2766.0002  # M132: -2316.9998 -> -2766.0002
$ws.Cells.Item(134, 8).Value = 49996.5  # H134: 53331 -> 49996.5
$ws.Cells.Item(134, 10).Value = 49996.5  # J134: 53331 -> 49996.5
$ws.Cells.Item(134, 12).Value = 149989.5  # L134: 159993 -> 149989.5
